$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin price / volume table cells to match the latest scrape.
# D-column cells whose new value looks numeric get an explicit Text
# number format first, so Excel does not silently reinterpret them as
# numbers (which would drop significant trailing zeros, e.g. "5.300").

$ws.Range("D2").Value = '30.694.34'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '1.911.61'
$ws.Range("E3").Value = '  +0.81%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.25'
$ws.Range("E5").Value = '  -1.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.004'
$ws.Range("E6").Value = '  +0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4926'
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  +1.20%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06740'
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = '1.922.06'
$ws.Range("E10").Value = '  +1.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '17.01'
$ws.Range("E11").Value = '  -1.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07354'
$ws.Range("E12").Value = '  +1.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.151'
$ws.Range("E13").Value = '  +2.42%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.12'
$ws.Range("E14").Value = '  -2.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6697'
$ws.Range("E15").Value = '  -0.85%  '
$ws.Range("D16").Value = '30.637.70'
$ws.Range("E16").Value = '  -0.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000007887'
$ws.Range("E17").Value = '  -1.08%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.46'
$ws.Range("E18").Value = '  +2.80%  '
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("D20").Value = '2.167.10'
$ws.Range("E20").Value = '  +1.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.300'
$ws.Range("E21").Value = '  +10.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.005'
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '193.74'
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.238'
$ws.Range("E24").Value = '  +2.60%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.637'
$ws.Range("E25").Value = '  +2.87%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '163.71'
$ws.Range("E26").Value = '  +4.75%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.51'
$ws.Range("E27").Value = '  -2.50%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.942'
$ws.Range("E28").Value = '  +2.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.491'
$ws.Range("E29").Value = '  +5.82%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.333'
$ws.Range("E30").Value = '  +0.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09102'
$ws.Range("E31").Value = '  +0.27%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.046'
$ws.Range("E32").Value = '  +1.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05227'
$ws.Range("E33").Value = '  +0.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7378'
$ws.Range("E34").Value = '  -0.41%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.108'
$ws.Range("E35").Value = '  +0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.729'
$ws.Range("E36").Value = '  -1.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01819'
$ws.Range("E37").Value = '  -0.73%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.715'
$ws.Range("E38").Value = '  +1.39%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9183'
$ws.Range("E39").Value = '  -1.19%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.055'
$ws.Range("E40").Value = '  -2.98%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.61'
$ws.Range("E41").Value = '  +29.63%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4430'
$ws.Range("E42").Value = '  +0.88%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '106.55'
$ws.Range("E43").Value = '  +1.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.895'
$ws.Range("E44").Value = '  +2.81%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.1382'
$ws.Range("E46").Value = '  +2.47%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.558'
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.057'
$ws.Range("E48").Value = '  +4.23%  '
$ws.Range("B49").Value = 'Elrond'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '35.36'
$ws.Range("E49").Value = '  +5.20%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05869'
$ws.Range("E50").Value = '  +0.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3998'
$ws.Range("E51").Value = '  +1.95%  '
